$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 8 (weekly re-sort by Fecha)
$ws.Range("D2").Value = 44511
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = 1350
$ws.Range("N2").Value = '$/kilo'
$ws.Range("P2").Value = 1350

# Row 3 <- source row 18 (weekly re-sort by Fecha)
$ws.Range("D3").Value = 44860
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1609
$ws.Range("P3").Value = 1609

# Row 4 <- source row 14 (weekly re-sort by Fecha)
$ws.Range("D4").Value = 44519
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1240
$ws.Range("P4").Value = 1240

# Row 5 <- source row 3 (weekly re-sort by Fecha)
$ws.Range("D5").Value = 44510
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1350
$ws.Range("P5").Value = 1350

# Row 6 <- source row 12 (weekly re-sort by Fecha)
$ws.Range("D6").Value = 44477
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1400
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1460
$ws.Range("P6").Value = 1460

# Row 7 <- source row 5 (weekly re-sort by Fecha)
$ws.Range("D7").Value = 44881
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 2600
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = 2650
$ws.Range("P7").Value = 2650

# Row 8 <- source row 6 (weekly re-sort by Fecha)
$ws.Range("D8").Value = 44881
$ws.Range("I8").Value = 'Segunda'
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 2400
$ws.Range("L8").Value = 2400
$ws.Range("M8").Value = 2400
$ws.Range("P8").Value = 2400

# Row 9 <- source row 15 (weekly re-sort by Fecha)
$ws.Range("D9").Value = 44876
$ws.Range("J9").Value = 350
$ws.Range("M9").Value = 1557
$ws.Range("O9").Value = 'Provincia de Linares'
$ws.Range("P9").Value = 1557

# Row 10 <- source row 13 (weekly re-sort by Fecha)
$ws.Range("D10").Value = 44468
$ws.Range("H10").Value = 'Verde'
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 1800
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1920
$ws.Range("P10").Value = 1920

# Row 12 <- source row 2 (weekly re-sort by Fecha)
$ws.Range("D12").Value = 44496
$ws.Range("J12").Value = 550
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 1773
$ws.Range("N12").Value = '$/paquete'
$ws.Range("P12").Value = 1773

# Row 13 <- source row 9 (weekly re-sort by Fecha)
$ws.Range("D13").Value = 44524
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1550
$ws.Range("O13").Value = 'Provincia de Talca'
$ws.Range("P13").Value = 1550

# Row 14 <- source row 10 (weekly re-sort by Fecha)
$ws.Range("D14").Value = 44489
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 1400
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = 1450
$ws.Range("P14").Value = 1450

# Row 15 <- source row 19 (weekly re-sort by Fecha)
$ws.Range("D15").Value = 44839
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 1700
$ws.Range("L15").Value = 1800
$ws.Range("M15").Value = 1760
$ws.Range("P15").Value = 1760

# Row 18 <- source row 4 (weekly re-sort by Fecha)
$ws.Range("D18").Value = 44875
$ws.Range("J18").Value = 300
$ws.Range("L18").Value = 1600
$ws.Range("M18").Value = 1550
$ws.Range("P18").Value = 1550

# Row 19 <- source row 7 (weekly re-sort by Fecha)
$ws.Range("D19").Value = 44526
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = 1550
$ws.Range("P19").Value = 1550

